$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$shp = $hdr.Range.InlineShapes.Item(1)
Write-Output ("HasChart=" + $shp.HasChart)
Write-Output ("Type=" + $shp.Type)
try { Write-Output ("Field=" + $shp.Field) } catch { Write-Output "Field error" }
try { Write-Output ("LinkFormat=" + $shp.LinkFormat) } catch { Write-Output "LinkFormat error" }
